$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Shift rows 10-46 up to rows 4-40 (delete old rows 4-9 worth of data, preserving cell types/styles via Copy)
for ($r = 4; $r -le 40; $r++) {
    $srcRow = $r + 6
    $src = $ws.Range("B" + $srcRow + ":F" + $srcRow)
    $dst = $ws.Range("B" + $r + ":F" + $r)
    $src.Copy($dst)
}

# Populate new rows 41-46 with the newly added translation entries
# Row 41
$ws.Cells.Item(41, 2).Value = "SingleUseId51"
$ws.Cells.Item(41, 3).Value = "Large"
$ws.Cells.Item(41, 4).Value = "Right"
$ws.Cells.Item(41, 5).Value = "<value> "
$ws.Cells.Item(41, 6).Value = "LTR"

# Row 42
$ws.Cells.Item(42, 2).Value = "SingleUseId52"
$ws.Cells.Item(42, 3).Value = "Large"
$ws.Cells.Item(42, 4).Value = "Left"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "125"
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(42, 6).Value = "LTR"

# Row 43
$ws.Cells.Item(43, 2).Value = "SingleUseId53"
$ws.Cells.Item(43, 3).Value = "Values"
$ws.Cells.Item(43, 4).Value = "Center"
$ws.Cells.Item(43, 5).Value = "<value>"
$ws.Cells.Item(43, 6).Value = "LTR"

# Row 44
$ws.Cells.Item(44, 2).Value = "SingleUseId54"
$ws.Cells.Item(44, 3).Value = "Values"
$ws.Cells.Item(44, 4).Value = "Left"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "120"
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(44, 6).Value = "LTR"

# Row 45
$ws.Cells.Item(45, 2).Value = "SingleUseId55"
$ws.Cells.Item(45, 3).Value = "Values"
$ws.Cells.Item(45, 4).Value = "Center"
$ws.Cells.Item(45, 5).Value = "<value>"
$ws.Cells.Item(45, 6).Value = "LTR"

# Row 46
$ws.Cells.Item(46, 2).Value = "SingleUseId56"
$ws.Cells.Item(46, 3).Value = "Values"
$ws.Cells.Item(46, 4).Value = "Left"
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "250"
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(46, 6).Value = "LTR"

Write-Output "Edit applied successfully"